# Update "想去人数" (want-to-go count) values in the F column for a handful
# of rows on both the "展览" sheet and the "全部类型" sheet, matching the
# upstream data refresh captured in the commit's XML diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value  = 101
$wsExhibit.Range("F8").Value  = 1768
$wsExhibit.Range("F12").Value = 2022
$wsExhibit.Range("F13").Value = 15
$wsExhibit.Range("F14").Value = 140
$wsExhibit.Range("F15").Value = 1225
$wsExhibit.Range("F16").Value = 459
$wsExhibit.Range("F19").Value = 208
$wsExhibit.Range("F23").Value = 46
$wsExhibit.Range("F25").Value = 7
$wsExhibit.Range("F26").Value = 1093

# Sheet "全部类型" (rId4 / sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 101
$wsAll.Range("F8").Value  = 1768
$wsAll.Range("F13").Value = 2022
$wsAll.Range("F14").Value = 15
$wsAll.Range("F15").Value = 140
$wsAll.Range("F16").Value = 1225
$wsAll.Range("F17").Value = 459
$wsAll.Range("F20").Value = 208
$wsAll.Range("F24").Value = 46
$wsAll.Range("F26").Value = 7
$wsAll.Range("F27").Value = 1093
